# Refresh the cryptocurrency price / 1h-volume columns.
# Data pulled from coinranking.com on Thu Apr 27 17:54:52 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Mimic typing a leading apostrophe in Excel so the numeric-
    # looking price strings (e.g. "0.9988") are stored as TEXT,
    # not auto-converted to the Number type.
    $ws.Range($addr).Value = "'" + $text
}

# Row 2
Set-TextCell "D2" '29.338.70'
$ws.Range("E2").Value = '  -1.46%  '

# Row 3
Set-TextCell "D3" '1.906.92'
$ws.Range("E3").Value = '  -2.34%  '

# Row 4
Set-TextCell "D4" '0.9988'
$ws.Range("E4").Value = '  -0.30%  '

# Row 5
Set-TextCell "D5" '333.86'
$ws.Range("E5").Value = '  -2.25%  '

# Row 6
Set-TextCell "D6" '0.9983'
$ws.Range("E6").Value = '  -0.32%  '

# Row 7
Set-TextCell "D7" '0.4639'
$ws.Range("E7").Value = '  -2.86%  '

# Row 8
Set-TextCell "D8" '0.4143'
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
Set-TextCell "D9" '47.89'
$ws.Range("E9").Value = '  +0.21%  '

# Row 10
Set-TextCell "D10" '0.08050'
$ws.Range("E10").Value = '  -2.34%  '

# Row 11
Set-TextCell "D11" '1.022'
$ws.Range("E11").Value = '  -1.40%  '

# Row 12
Set-TextCell "D12" '22.29'
$ws.Range("E12").Value = '  -2.00%  '

# Row 13
Set-TextCell "D13" '1.927.91'
$ws.Range("E13").Value = '  -1.28%  '

# Row 14
Set-TextCell "D14" '5.969'
$ws.Range("E14").Value = '  -3.23%  '

# Row 15
Set-TextCell "D15" '7.147'
$ws.Range("E15").Value = '  -3.32%  '

# Row 16
Set-TextCell "D16" '89.28'
$ws.Range("E16").Value = '  -2.95%  '

# Row 17
Set-TextCell "D17" '0.9992'
$ws.Range("E17").Value = '  -0.38%  '

# Row 18
Set-TextCell "D18" '0.00001034'
$ws.Range("E18").Value = '  -2.27%  '

# Row 19
Set-TextCell "D19" '0.06589'
$ws.Range("E19").Value = '  -1.62%  '

# Row 20
Set-TextCell "D20" '17.71'
$ws.Range("E20").Value = '  -1.74%  '

# Row 21
Set-TextCell "D21" '0.9993'
$ws.Range("E21").Value = '  -0.11%  '

# Row 22
Set-TextCell "D22" '29.230.42'
$ws.Range("E22").Value = '  -1.76%  '

# Row 23
Set-TextCell "D23" '5.513'
$ws.Range("E23").Value = '  -1.34%  '

# Row 24
Set-TextCell "D24" '11.43'
$ws.Range("E24").Value = '  +1.46%  '

# Row 25
$ws.Range("E25").Value = '  -3.62%  '

# Row 26
Set-TextCell "D26" '2.171.36'
$ws.Range("E26").Value = '  -0.42%  '

# Row 27
Set-TextCell "D27" '156.68'
$ws.Range("E27").Value = '  -2.75%  '

# Row 28
Set-TextCell "D28" '19.88'
$ws.Range("E28").Value = '  -1.77%  '

# Row 29
Set-TextCell "D29" '2.148'
$ws.Range("E29").Value = '  -1.39%  '

# Row 30
Set-TextCell "D30" '5.669'
$ws.Range("E30").Value = '  -0.01%  '

# Row 31
Set-TextCell "D31" '117.45'
$ws.Range("E31").Value = '  -4.50%  '

# Row 32
Set-TextCell "D32" '1.042'
$ws.Range("E32").Value = '  +3.83%  '

# Row 33
Set-TextCell "D33" '0.09460'
$ws.Range("E33").Value = '  -1.92%  '

# Row 34
Set-TextCell "D34" '1.431'
$ws.Range("E34").Value = '  -2.98%  '

# Row 35
Set-TextCell "D35" '3.546'
$ws.Range("E35").Value = '  -3.83%  '

# Row 36
Set-TextCell "D36" '5.398'
$ws.Range("E36").Value = '  -1.89%  '

# Row 37
Set-TextCell "D37" '0.06116'
$ws.Range("E37").Value = '  -2.41%  '

# Row 38
Set-TextCell "D38" '0.02262'
$ws.Range("E38").Value = '  -2.21%  '

# Row 39
$ws.Range("E39").Value = '  +0.02%  '

# Row 40
$ws.Range("E40").Value = '  -0.17%  '

# Row 41
Set-TextCell "D41" '0.5893'
$ws.Range("E41").Value = '  -3.21%  '

# Row 42
Set-TextCell "D42" '0.9982'
$ws.Range("E42").Value = '  -0.27%  '

# Row 43
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell "D43" '10.21'
$ws.Range("E43").Value = '  -4.89%  '

# Row 44
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell "D44" '0.1833'
$ws.Range("E44").Value = '  -3.10%  '

# Row 45
Set-TextCell "D45" '2.378'
$ws.Range("E45").Value = '  -0.09%  '

# Row 46
Set-TextCell "D46" '1.255'
$ws.Range("E46").Value = '  -1.89%  '

# Row 47
Set-TextCell "D47" '0.07522'
$ws.Range("E47").Value = '  +1.74%  '

# Row 48
Set-TextCell "D48" '0.5576'
$ws.Range("E48").Value = '  -2.38%  '

# Row 49
Set-TextCell "D49" '12.17'

# Row 50
Set-TextCell "D50" '1.933'
$ws.Range("E50").Value = '  -2.71%  '

# Row 51
Set-TextCell "D51" '112.85'
$ws.Range("E51").Value = '  -0.16%  '
